$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.402.12"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.346.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.338.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.55"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.883.67"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.343.41"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.407.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.978"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.01"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.24%  "
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.23"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -4.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.72"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "577.46"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.39"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.24%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -8.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.66"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("E41").Value = "  -4.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.100.34"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  -5.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.23"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0409"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  -3.41%  "
